{"js": "// Replacement pairs extracted from the diff, in document order\n// (the document's paragraph order: index 0 is the date line, the\n// remaining paragraphs are the 100 arithmetic expressions in the\n// table cells, in row-major order matching these pairs).\nconst replacements = [\n  [\"50-4=\", \"8+41=\"],\n  [\"5+77=\", \"89-59=\"],\n  [\"71+11=\", \"30+28=\"],\n  [\"65-54=\", \"17+3=\"],\n  [\"85+5=\", \"98-22=\"],\n  [\"22+8=\", \"27+12=\"],\n  [\"62+14=\", \"16+2=\"],\n  [\"33+26=\", \"44-23=\"],\n  [\"92-15=\", \"9+78=\"],\n  [\"36-19=\", \"47+26=\"],\n  [\"80-34=\", \"27+53=\"],\n  [\"72-41=\", \"61-22=\"],\n  [\"94-69=\", \"24+59=\"],\n  [\"69-50=\", \"90-76=\"],\n  [\"50-35=\", \"12+16=\"],\n  [\"70-46=\", \"5+67=\"],\n  [\"42+10=\", \"64+26=\"],\n  [\"10+59=\", \"80-54=\"],\n  [\"9+79=\", \"87-28=\"],\n  [\"3+13=\", \"7+3=\"],\n  [\"71-28=\", \"37+58=\"],\n  [\"53-17=\", \"19+37=\"],\n  [\"96-32=\", \"35+34=\"],\n  [\"14+31=\", \"44-17=\"],\n  [\"37+26=\", \"13-6=\"],\n  [\"29-14=\", \"62-45=\"],\n  [\"50+8=\", \"81+0=\"],\n  [\"15+44=\", \"76-71=\"],\n  [\"48+25=\", \"21+11=\"],\n  [\"43+44=\", \"30+21=\"],\n  [\"14+84=\", \"50+21=\"],\n  [\"18+68=\", \"11+84=\"],\n  [\"44-32=\", \"56+13=\"],\n  [\"56-4=\", \"41-11=\"],\n  [\"68-43=\", \"9+44=\"],\n  [\"94-89=\", \"17+70=\"],\n  [\"24+49=\", \"5+2=\"],\n  [\"73+10=\", \"77+19=\"],\n  [\"40-21=\", \"49+35=\"],\n  [\"90-50=\", \"95-62=\"],\n  [\"85-43=\", \"19+0=\"],\n  [\"12+52=\", \"44-42=\"],\n  [\"58-3=\", \"51-8=\"],\n  [\"68+30=\", \"76-15=\"],\n  [\"41-20=\", \"80-71=\"],\n  [\"74-12=\", \"76-31=\"],\n  [\"1+87=\", \"82-81=\"],\n  [\"28-3=\", \"53+12=\"],\n  [\"27-11=\", \"65-11=\"],\n  [\"72-7=\", \"64+28=\"],\n  [\"87-12=\", \"30-8=\"],\n  [\"67+24=\", \"20-16=\"],\n  [\"20+71=\", \"2+86=\"],\n  [\"88-11=\", \"15-14=\"],\n  [\"85-33=\", \"56+22=\"],\n  [\"74+21=\", \"61-34=\"],\n  [\"4+52=\", \"65+31=\"],\n  [\"4+95=\", \"60+20=\"],\n  [\"72-32=\", \"81-77=\"],\n  [\"58-0=\", \"47-18=\"],\n  [\"16+10=\", \"37+46=\"],\n  [\"40+55=\", \"54-12=\"],\n  [\"28+45=\", \"73-39=\"],\n  [\"88+8=\", \"21+69=\"],\n  [\"70-1=\", \"81-30=\"],\n  [\"34-24=\", \"27+17=\"],\n  [\"21-9=\", \"12+73=\"],\n  [\"44+4=\", \"35-29=\"],\n  [\"89-69=\", \"86+1=\"],\n  [\"95-16=\", \"87-48=\"],\n  [\"25-17=\", \"32-28=\"],\n  [\"71-47=\", \"4+27=\"],\n  [\"59+30=\", \"98-0=\"],\n  [\"15+33=\", \"30+6=\"],\n  [\"49-15=\", \"17+70=\"],\n  [\"21+24=\", \"25-23=\"],\n  [\"97-13=\", \"22+3=\"],\n  [\"39-0=\", \"48+32=\"],\n  [\"50-6=\", \"56+16=\"],\n  [\"70-32=\", \"40-3=\"],\n  [\"96-21=\", \"25+38=\"],\n  [\"96-3=\", \"50+32=\"],\n  [\"99-82=\", \"74-29=\"],\n  [\"12+50=\", \"29+65=\"],\n  [\"9+24=\", \"26+14=\"],\n  [\"83-38=\", \"23+4=\"],\n  [\"68+6=\", \"66+25=\"],\n  [\"80-45=\", \"0+5=\"],\n  [\"25+68=\", \"5+76=\"],\n  [\"25+37=\", \"18+11=\"],\n  [\"93-10=\", \"39-34=\"],\n  [\"93-13=\", \"24+60=\"],\n  [\"38+18=\", \"41-35=\"],\n  [\"34+33=\", \"67-9=\"],\n  [\"45+29=\", \"11+34=\"],\n  [\"1+12=\", \"33+56=\"],\n  [\"15+55=\", \"64+11=\"],\n  [\"42-25=\", \"34+55=\"],\n  [\"71+16=\", \"95-20=\"],\n  [\"44+23=\", \"63-56=\"]\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// The first paragraph is the date line (\"2024-10-11 Friday\") and is\n// left untouched. The remaining paragraphs correspond 1:1, in order,\n// to the `replacements` pairs above.\nconst exprParagraphs = paragraphs.items.slice(1);\n\nif (exprParagraphs.length !== replacements.length) {\n  throw new Error(\n    \"Unexpected paragraph count: expected \" + replacements.length +\n    \" expression paragraphs, found \" + exprParagraphs.length\n  );\n}\n\nfor (let i = 0; i < replacements.length; i++) {\n  const [oldText, newText] = replacements[i];\n  const para = exprParagraphs[i];\n  if (para.text !== oldText) {\n    throw new Error(\n      \"Paragraph \" + i + \" text mismatch: expected '\" + oldText +\n      \"' but found '\" + para.text + \"'\"\n    );\n  }\n  para.insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Replacement pairs extracted from the diff, in document order.\n# Each pair is a unique \"old\" expression text (e.g. \"50-4=\") mapped to\n# its replacement (e.g. \"8+41=\"), applied as a whole-document\n# Find & Replace so the run-level formatting (font/size) of each cell\n# is preserved untouched.\n$replacements = @(\n    @{Old=\"50-4=\"; New=\"8+41=\"},\n    @{Old=\"5+77=\"; New=\"89-59=\"},\n    @{Old=\"71+11=\"; New=\"30+28=\"},\n    @{Old=\"65-54=\"; New=\"17+3=\"},\n    @{Old=\"85+5=\"; New=\"98-22=\"},\n    @{Old=\"22+8=\"; New=\"27+12=\"},\n    @{Old=\"62+14=\"; New=\"16+2=\"},\n    @{Old=\"33+26=\"; New=\"44-23=\"},\n    @{Old=\"92-15=\"; New=\"9+78=\"},\n    @{Old=\"36-19=\"; New=\"47+26=\"},\n    @{Old=\"80-34=\"; New=\"27+53=\"},\n    @{Old=\"72-41=\"; New=\"61-22=\"},\n    @{Old=\"94-69=\"; New=\"24+59=\"},\n    @{Old=\"69-50=\"; New=\"90-76=\"},\n    @{Old=\"50-35=\"; New=\"12+16=\"},\n    @{Old=\"70-46=\"; New=\"5+67=\"},\n    @{Old=\"42+10=\"; New=\"64+26=\"},\n    @{Old=\"10+59=\"; New=\"80-54=\"},\n    @{Old=\"9+79=\"; New=\"87-28=\"},\n    @{Old=\"3+13=\"; New=\"7+3=\"},\n    @{Old=\"71-28=\"; New=\"37+58=\"},\n    @{Old=\"53-17=\"; New=\"19+37=\"},\n    @{Old=\"96-32=\"; New=\"35+34=\"},\n    @{Old=\"14+31=\"; New=\"44-17=\"},\n    @{Old=\"37+26=\"; New=\"13-6=\"},\n    @{Old=\"29-14=\"; New=\"62-45=\"},\n    @{Old=\"50+8=\"; New=\"81+0=\"},\n    @{Old=\"15+44=\"; New=\"76-71=\"},\n    @{Old=\"48+25=\"; New=\"21+11=\"},\n    @{Old=\"43+44=\"; New=\"30+21=\"},\n    @{Old=\"14+84=\"; New=\"50+21=\"},\n    @{Old=\"18+68=\"; New=\"11+84=\"},\n    @{Old=\"44-32=\"; New=\"56+13=\"},\n    @{Old=\"56-4=\"; New=\"41-11=\"},\n    @{Old=\"68-43=\"; New=\"9+44=\"},\n    @{Old=\"94-89=\"; New=\"17+70=\"},\n    @{Old=\"24+49=\"; New=\"5+2=\"},\n    @{Old=\"73+10=\"; New=\"77+19=\"},\n    @{Old=\"40-21=\"; New=\"49+35=\"},\n    @{Old=\"90-50=\"; New=\"95-62=\"},\n    @{Old=\"85-43=\"; New=\"19+0=\"},\n    @{Old=\"12+52=\"; New=\"44-42=\"},\n    @{Old=\"58-3=\"; New=\"51-8=\"},\n    @{Old=\"68+30=\"; New=\"76-15=\"},\n    @{Old=\"41-20=\"; New=\"80-71=\"},\n    @{Old=\"74-12=\"; New=\"76-31=\"},\n    @{Old=\"1+87=\"; New=\"82-81=\"},\n    @{Old=\"28-3=\"; New=\"53+12=\"},\n    @{Old=\"27-11=\"; New=\"65-11=\"},\n    @{Old=\"72-7=\"; New=\"64+28=\"},\n    @{Old=\"87-12=\"; New=\"30-8=\"},\n    @{Old=\"67+24=\"; New=\"20-16=\"},\n    @{Old=\"20+71=\"; New=\"2+86=\"},\n    @{Old=\"88-11=\"; New=\"15-14=\"},\n    @{Old=\"85-33=\"; New=\"56+22=\"},\n    @{Old=\"74+21=\"; New=\"61-34=\"},\n    @{Old=\"4+52=\"; New=\"65+31=\"},\n    @{Old=\"4+95=\"; New=\"60+20=\"},\n    @{Old=\"72-32=\"; New=\"81-77=\"},\n    @{Old=\"58-0=\"; New=\"47-18=\"},\n    @{Old=\"16+10=\"; New=\"37+46=\"},\n    @{Old=\"40+55=\"; New=\"54-12=\"},\n    @{Old=\"28+45=\"; New=\"73-39=\"},\n    @{Old=\"88+8=\"; New=\"21+69=\"},\n    @{Old=\"70-1=\"; New=\"81-30=\"},\n    @{Old=\"34-24=\"; New=\"27+17=\"},\n    @{Old=\"21-9=\"; New=\"12+73=\"},\n    @{Old=\"44+4=\"; New=\"35-29=\"},\n    @{Old=\"89-69=\"; New=\"86+1=\"},\n    @{Old=\"95-16=\"; New=\"87-48=\"},\n    @{Old=\"25-17=\"; New=\"32-28=\"},\n    @{Old=\"71-47=\"; New=\"4+27=\"},\n    @{Old=\"59+30=\"; New=\"98-0=\"},\n    @{Old=\"15+33=\"; New=\"30+6=\"},\n    @{Old=\"49-15=\"; New=\"17+70=\"},\n    @{Old=\"21+24=\"; New=\"25-23=\"},\n    @{Old=\"97-13=\"; New=\"22+3=\"},\n    @{Old=\"39-0=\"; New=\"48+32=\"},\n    @{Old=\"50-6=\"; New=\"56+16=\"},\n    @{Old=\"70-32=\"; New=\"40-3=\"},\n    @{Old=\"96-21=\"; New=\"25+38=\"},\n    @{Old=\"96-3=\"; New=\"50+32=\"},\n    @{Old=\"99-82=\"; New=\"74-29=\"},\n    @{Old=\"12+50=\"; New=\"29+65=\"},\n    @{Old=\"9+24=\"; New=\"26+14=\"},\n    @{Old=\"83-38=\"; New=\"23+4=\"},\n    @{Old=\"68+6=\"; New=\"66+25=\"},\n    @{Old=\"80-45=\"; New=\"0+5=\"},\n    @{Old=\"25+68=\"; New=\"5+76=\"},\n    @{Old=\"25+37=\"; New=\"18+11=\"},\n    @{Old=\"93-10=\"; New=\"39-34=\"},\n    @{Old=\"93-13=\"; New=\"24+60=\"},\n    @{Old=\"38+18=\"; New=\"41-35=\"},\n    @{Old=\"34+33=\"; New=\"67-9=\"},\n    @{Old=\"45+29=\"; New=\"11+34=\"},\n    @{Old=\"1+12=\"; New=\"33+56=\"},\n    @{Old=\"15+55=\"; New=\"64+11=\"},\n    @{Old=\"42-25=\"; New=\"34+55=\"},\n    @{Old=\"71+16=\"; New=\"95-20=\"},\n    @{Old=\"44+23=\"; New=\"63-56=\"}\n)\n\n$d = $word.ActiveDocument\n\n$replacedCount = 0\n\nforeach ($pair in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $found = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n    if ($found) {\n        $replacedCount = $replacedCount + 1\n    } else {\n        Write-Output \"WARNING: replacement not found for '$($pair.Old)'\"\n    }\n}\n\nWrite-Output \"Replaced $replacedCount of $($replacements.Count) expressions\"\n"}
